$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row before row 23 ---
# This pushes the (empty) rows 23-26 and the signature block (rows 27-28)
# down by one, becoming rows 24-27 and 28-29 respectively, and Excel
# automatically re-points the existing merged cells.
$ws.Rows.Item(23).Insert()

# --- 2. Give the new row 23 the "last data row" look (thicker bottom
#        border) that row 22 used to have ---
$ws.Range("B22:J22").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 3. Row 22 becomes a normal/middle data row again (like row 21) ---
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 4. Fill the new row 23 with the additional period (2509) for the
#        same worker as row 22 (OSCAR DAVID ANILLO GUZMAN) ---
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1051818008"
$ws.Range("D23").Value = "OSCAR DAVID ANILLO GUZMAN"
$ws.Range("E23").Value = "2509"
$ws.Range("F23").Value = 60000
$ws.Range("G23").Value = 1500000

# --- 5. Center-align the "Periodo Mora" column (E) across every data
#        row, including the newly added one ---
$ws.Range("E16:E23").HorizontalAlignment = -4108   # xlCenter

# --- 6. Update the totals on the summary block ---
$ws.Range("E11").Value = 372500   # VALOR MORA
$ws.Range("F13").Value = 6        # Cant. Periodos
